$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.756.37'
$ws.Range("E2").Value = '  +1.93%  '
$ws.Range("D3").Value = '1.575.04'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("D5").Value = '213.24'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '44.98'
$ws.Range("E8").Value = '  +2.49%  '
$ws.Range("D9").Value = '24.22'
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("E10").Value = '  -1.24%  '
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("D12").Value = '0.0888'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '1.798.82'
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").Value = '1.565.92'
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '28.731.74'
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").Value = '0.522'
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '3.69'
$ws.Range("E17").Value = '  -1.56%  '
$ws.Range("D18").Value = "'62.50"
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").Value = '230.41'
$ws.Range("D21").Value = '0.0₃0694'
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D24").Value = '9.19'
$ws.Range("E25").Value = '  +8.55%  '
$ws.Range("D26").Value = '151.84'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("E29").Value = '  -2.56%  '
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").Value = '0.0485'
$ws.Range("E31").Value = '  +2.62%  '
$ws.Range("E32").Value = '  -1.95%  '
$ws.Range("D33").Value = '3.21'
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("D35").Value = '1.390.03'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("D36").Value = '1.05'
$ws.Range("E36").Value = '  +2.82%  '
$ws.Range("D37").Value = '1.55'
$ws.Range("E37").Value = '  -2.99%  '
$ws.Range("D38").Value = '2.36'
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("E39").Value = '  +2.85%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("E41").Value = '  -2.65%  '
$ws.Range("D42").Value = '1.91'
$ws.Range("E42").Value = '  +2.21%  '
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("D45").Value = "'0.0470"
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("E46").Value = '  -1.22%  '
$ws.Range("E47").Value = '  -1.36%  '
$ws.Range("D48").Value = '63.42'
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("D49").Value = '1.711.88'
$ws.Range("E49").Value = '  -0.78%  '
$ws.Range("D50").Value = '86.73'
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("E51").Value = '  -0.62%  '
